{"js": "// Remove the duplicated title-page / preface block that appears at the\n// very start of the document (12 paragraphs: 3 blank lines, \"2024\",\n// \"OFFICIAL PLAYING RULES\", \"OF THE\", \"WALLERSTEIN FANTASY FOOTBALL\n// LEAGUE\", the logo image + page break paragraph, \"PREFACE\", the two\n// preface text paragraphs, and the trailing page-break paragraph).\n// The document contained this block twice in a row; the fix drops the\n// first (incomplete) copy and keeps the second copy as the real start\n// of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst REMOVE_COUNT = 12;\nconst count = Math.min(REMOVE_COUNT, paragraphs.items.length);\nfor (let i = 0; i < count; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The document starts with the title-page / preface block (3 blank\n# centered lines, \"2024\", \"OFFICIAL PLAYING RULES\", \"OF THE\",\n# \"WALLERSTEIN FANTASY FOOTBALL LEAGUE\", the league-logo image with a\n# trailing page break, \"PREFACE\", the two preface paragraphs, and a\n# final page-break paragraph -- 12 paragraphs total) duplicated back to\n# back. The fix removes the first (incomplete/stray) copy of that\n# 12-paragraph block and keeps the second copy, which is followed\n# immediately by \"SECTION 1. KEEPERS\".\n\n$d = $word.ActiveDocument\n\n$removeCount = 12\nfor ($i = 0; $i -lt $removeCount; $i++) {\n    $d.Paragraphs(1).Range.Delete()\n}\n"}
